# Artfynd sheet: the three species-observation rows (2,3,4) got rotated
# ("row 3" data moved up into row 2, "row 4" data moved up into row 3,
# and the old "row 2" data wrapped around into row 4), and a previously
# missing "Samlings-nummer" (AR) value was filled in for what is now row 2.
#
# Rather than re-deriving that rotation generically, we just write the
# final literal values for every touched cell - this keeps numeric
# precision exact and lets us control text-vs-number typing per cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay TEXT even though it looks like a
# number/date (Excel would otherwise silently coerce it). Using a leading
# apostrophe forces text entry (exactly like typing '10 into a cell), then
# we restore the Normal style so no stray "quote prefix" formatting lingers.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---- Row 2 (now holds the former row-3 record) ----
$ws.Range("A2").Value = 66541020
$ws.Range("B2").Value = 89392
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = "Ullticka"
$ws.Range("G2").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P2").Value = "Grysjöbäcken, Mpd"
$ws.Range("Q2").Value = 633480.2438334802
$ws.Range("R2").Value = 6940243.949297423
$ws.Range("S2").Value = 50
Set-TextValue $ws.Range("Y2") "2014-09-18"
Set-TextValue $ws.Range("AA2") "2014-09-18"
Set-TextValue $ws.Range("AR2") "23349"
$ws.Range("AW2").Value = "Malin Sahlin"
$ws.Range("AX2").Value = "Via Malin Sahlin"

# ---- Row 3 (now holds the former row-4 record) ----
$ws.Range("A3").Value = 66541021
$ws.Range("B3").Value = 77506
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
Set-TextValue $ws.Range("I3") "10"
Set-TextValue $ws.Range("AR3") "23350"

# ---- Row 4 (now holds the former row-2 record) ----
$ws.Range("A4").Value = 91961108
$ws.Range("B4").Value = 78569
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
Set-TextValue $ws.Range("I4") ""
$ws.Range("P4").Value = "Mysjöberget, Mpd"
$ws.Range("Q4").Value = 632736.7829380766
$ws.Range("R4").Value = 6940262.09546657
$ws.Range("S4").Value = 25
Set-TextValue $ws.Range("Y4") "2020-06-09"
Set-TextValue $ws.Range("AA4") "2020-06-09"
$ws.Range("AR4").Value = ""
$ws.Range("AW4").Value = "Mikael Gudrunsson"
$ws.Range("AX4").Value = "Mikael Gudrunsson"
